# Update workload status on the GanttChart sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")

# Move the "Display Week" scroll bar back one week (5 -> 4). Updating the
# linked cell mirrors what moving the scrollbar control does, and also
# keeps the embedded form-control's cached value (ctrlProp1.xml) in sync.
$ws.Range("H4").Value = 4

# Mark a few tasks as progressed / completed (% DONE column).
$ws.Range("H21").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("H24").Value = 1

# Scroll the frozen pane down and move the active selection, matching the
# view state saved with the workbook.
$ws.Activate()
$ws.Range("A20").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("T30").Select()
